$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 23: Friday (H23) now marked done -> style 18 (green fill, Friday border),
#     keep its "Lön" text; bump worked days from 4 to 5.
$ws.Range("H4").Copy()
$ws.Range("H23").PasteSpecial(-4122)
$ws.Range("I23").Value = 5

# --- Row 24: Tue-Fri (E24:H24) now marked done -> style 17 (green fill); 4 days worked.
$ws.Range("E4").Copy()
$ws.Range("E24:H24").PasteSpecial(-4122)
$ws.Range("I24").Value = 4

# --- Row 25: Mon-Fri (D25:H25) now marked done -> style 17 (green fill); 5 days worked.
$ws.Range("E4").Copy()
$ws.Range("D25:H25").PasteSpecial(-4122)
$ws.Range("I25").Value = 5

# --- Row 26: Mon-Tue (D26:E26) now marked done -> style 17 (green fill); 2 days worked.
$ws.Range("E4").Copy()
$ws.Range("D26:E26").PasteSpecial(-4122)
$ws.Range("I26").Value = 2

# --- Move the active selection from D31 to D30.
$ws.Range("D30").Select() | Out-Null
